# Update convergence table with rounded computational time values,
# and fill in the missing Bound_cond value on row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the missing boundary condition label for row 12 (F12 was blank, should be "Dirichlet")
$ws.Range("F12").Value = "Dirichlet"

# Round the computational time values in column J (rows 2-14)
$ws.Range("J2").Value = 105.362
$ws.Range("J3").Value = 6.3
$ws.Range("J4").Value = 209.234
$ws.Range("J5").Value = 11.911
$ws.Range("J6").Value = 9.772
$ws.Range("J7").Value = 9.846
$ws.Range("J8").Value = 15.677
$ws.Range("J9").Value = 15.956
$ws.Range("J10").Value = 4.784
$ws.Range("J11").Value = 2.487
$ws.Range("J12").Value = 5.976
$ws.Range("J13").Value = 62.252
$ws.Range("J14").Value = 3.775
